$d = $word.ActiveDocument

# The document currently contains three repeated copies of the same
# "Commit the code..." / "Push the code..." instruction pair (each
# followed by a blank paragraph), e.g.:
#   1: Commit the code by running the command: git commit ...
#   2: Push the code to the remote repo by running the command git push ...
#   3: (blank)
#   4: Commit the code by running the command: git commit ...
#   5: Push the code to the remote repo by running the command git push ...
#   6: (blank)
#   7: Commit the code by running the command: git commit ...
#   8: Push the code to the remote repo by running the command git push ...
#   9: (blank, holds the _GoBack bookmark)
#  10: (blank, trailing paragraph before the section break)
#
# Only the first pair should remain. Remove the duplicated 2nd and 3rd
# pairs (paragraphs 3 through 8 inclusive), keeping the bookmark
# paragraph (9) and the trailing blank paragraph (10) intact.
$pStart = $d.Paragraphs.Item(3)
$pEnd = $d.Paragraphs.Item(8)
$dupRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$dupRange.Delete()

# Each remaining instruction paragraph is split across three runs
# (e.g. "Commit the code...: " / "git" / " commit -m ...") with
# w:proofErr spell-check markers wrapped around the word "git". Use
# Find/Replace to re-type the full paragraph text in place, which
# collapses it down to a single run and drops the proofErr markers,
# matching the target formatting.
$find = $d.Content.Find
$find.ClearFormatting() | Out-Null
$find.Execute( `
    "Commit the code by running the command: git commit –m “committing the local changes”", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Commit the code by running the command: git commit –m “committing the local changes”", `
    2) | Out-Null

$find2 = $d.Content.Find
$find2.ClearFormatting() | Out-Null
$find2.Execute( `
    "Push the code to the remote repo by running the command git push origin master", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Push the code to the remote repo by running the command git push origin master", `
    2) | Out-Null
